$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '29.688.18'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  -3.15%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '2.097.05'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  -1.20%  '
$ws.Range("E4").Value = '  -0.54%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '343.16'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -2.82%  '
$ws.Range("E6").Value = '  -0.48%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.5138'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  -2.73%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.4400'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  -2.95%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '53.15'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  -1.49%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.09162'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +0.84%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '1.169'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -1.03%  '
$ws.Range("E12").Value = '  +1.35%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '2.100.86'
$ws.Range("D13").Style = "Normal"
$ws.Range("E14").Value = '  -1.45%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '8.164'
$ws.Range("D15").Style = "Normal"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '99.27'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  -3.07%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '1.008'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  -0.52%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '20.83'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +7.16%  '
$ws.Range("E20").Value = '  -1.19%  '
$ws.Range("E21").Value = '  -0.43%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '6.178'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -2.55%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '29.748.32'
$ws.Range("D23").Style = "Normal"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '12.56'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -2.00%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.289'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -4.21%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '2.347.30'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -0.98%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '21.79'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -2.86%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '162.57'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -1.67%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '2.508'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  -2.42%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '132.35'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  -2.97%  '
$ws.Range("E31").Value = '  -5.54%  '
$ws.Range("E32").Value = '  -2.98%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '1.634'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  -1.25%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '6.149'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -3.46%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '3.970'
$ws.Range("D35").Style = "Normal"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '6.063'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +0.83%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '10.15'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -1.82%  '
$ws.Range("E38").Value = '  -3.75%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.06695'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -2.81%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '12.38'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -1.36%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.6846'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -1.19%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.2212'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -4.69%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '1.296'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +1.31%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.6644'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +2.65%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '14.17'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  -3.98%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '2.294'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -1.73%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '3.610'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -4.29%  '
$ws.Range("B48").Value = 'BabyDogeCoin'
$ws.Range("C48").Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.00000000342'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -6.85%  '
$ws.Range("B49").Value = 'EOS'
$ws.Range("C49").Value = 'https://coinranking.com/coin/iAzbfXiBBKkR6+eos-eos'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '1.217'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -3.19%  '
$ws.Range("B50").Value = 'Aave'
$ws.Range("C50").Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '81.81'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  -1.33%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '1.159'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  -2.92%  '
